$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment is safe ---
$ws.Range("D2").Value = "70.580.94"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "2.523.46"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "2.523.24"
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").Value = "2.984.91"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "70.489.61"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("E17").Value = "  -5.19%  "
$ws.Range("D18").Value = "2.523.10"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("E19").Value = "  -5.72%  "
$ws.Range("E20").Value = "  -8.05%  "
$ws.Range("E21").Value = "  -4.24%  "
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("E28").Value = "  -5.34%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -5.68%  "
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("E42").Value = "  -6.53%  "
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("E44").Value = "  -5.62%  "
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("E47").Value = "  -9.42%  "
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("E49").Value = "  -6.33%  "
$ws.Range("E50").Value = "  -6.78%  "
$ws.Range("E51").Value = "  -1.26%  "

# --- Numeric-looking text values: must be forced to Text via a scratch cell ---
# (A1 is inside the existing used range, so round-tripping through it does not
#  trigger a structural dimension change; Clear() removes the residual Text
#  number-format picked up from the apostrophe-forced entry.)
$ws.Range("A1").Value = "'575.85"
$ws.Range("A1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'169.72"
$ws.Range("A1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'0.0000180"
$ws.Range("A1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'24.82"
$ws.Range("A1").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'11.49"
$ws.Range("A1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'7.53"
$ws.Range("A1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'355.47"
$ws.Range("A1").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'1.96"
$ws.Range("A1").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'69.19"
$ws.Range("A1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'9.19"
$ws.Range("A1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'1.00"
$ws.Range("A1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'478.13"
$ws.Range("A1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'1.28"
$ws.Range("A1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'157.55"
$ws.Range("A1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'18.85"
$ws.Range("A1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'18.57"
$ws.Range("A1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'1.30"
$ws.Range("A1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'4.70"
$ws.Range("A1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'38.28"
$ws.Range("A1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'141.24"
$ws.Range("A1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'0.523"
$ws.Range("A1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("A1").Value = "'0.595"
$ws.Range("A1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("A1").Clear()
